$d = $word.ActiveDocument

# Target the last paragraph in the document (the trailing empty paragraph
# after the "Live Coder" title) and replace its content with the new
# heading + body text, preserving its own paragraph mark/formatting by
# re-supplying equivalent <w:pPr>/<w:rPr> for every inserted paragraph.
$target = $d.Paragraphs($d.Paragraphs.Count)
$r = $target.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="sv-SE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>Förutsättningar</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve">Meningen är att låta programmet kommunicera med en underliggande kompilerare. På det sättet räcker det att skriva kod i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t>appen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> och sen när man sparar, som tex. i en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t>node</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t>app</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve">, så ska </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t>appen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> laddas om och man ser resultatet i ett fönster bredvid självaste live editorn.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sv-SE"/></w:rPr></w:pPr></w:p>
"@

$r.InsertXML($xml)
